$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells contain numeric-looking text that must remain text.
# Force text number-format before assignment, then restore default style so the
# cell keeps matching the original (un-styled) cells.
$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "30.260.02"
$rng.Style = "Normal"
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "1.887.79"
$rng.Style = "Normal"
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "238.34"
$rng.Style = "Normal"
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "0.4677"
$rng.Style = "Normal"
$rng = $ws.Range("D9")
$rng.NumberFormat = "@"
$rng.Value = "0.06605"
$rng.Style = "Normal"
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "20.03"
$rng.Style = "Normal"
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "0.07790"
$rng.Style = "Normal"
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "98.10"
$rng.Style = "Normal"
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "1.891.19"
$rng.Style = "Normal"
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "5.121"
$rng.Style = "Normal"
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "0.6789"
$rng.Style = "Normal"
$rng = $ws.Range("D16")
$rng.NumberFormat = "@"
$rng.Value = "285.88"
$rng.Style = "Normal"
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "30.278.13"
$rng.Style = "Normal"
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "1.001"
$rng.Style = "Normal"
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "12.66"
$rng.Style = "Normal"
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "2.142.49"
$rng.Style = "Normal"
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "5.395"
$rng.Style = "Normal"
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "0.000007314"
$rng.Style = "Normal"
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "6.193"
$rng.Style = "Normal"
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "9.432"
$rng.Style = "Normal"
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "166.25"
$rng.Style = "Normal"
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "19.29"
$rng.Style = "Normal"
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "2.000"
$rng.Style = "Normal"
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "1.376"
$rng.Style = "Normal"
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "0.09737"
$rng.Style = "Normal"
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "4.431"
$rng.Style = "Normal"
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "1.488"
$rng.Style = "Normal"
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "4.169"
$rng.Style = "Normal"
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "0.04686"
$rng.Style = "Normal"
$rng = $ws.Range("D35")
$rng.NumberFormat = "@"
$rng.Value = "0.7113"
$rng.Style = "Normal"
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = "1.099"
$rng.Style = "Normal"
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "2.712"
$rng.Style = "Normal"
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "0.01878"
$rng.Style = "Normal"
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "6.711"
$rng.Style = "Normal"
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "2.530"
$rng.Style = "Normal"
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "72.62"
$rng.Style = "Normal"
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "1.980"
$rng.Style = "Normal"
$rng = $ws.Range("D43")
$rng.NumberFormat = "@"
$rng.Value = "0.8688"
$rng.Style = "Normal"
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "104.34"
$rng.Style = "Normal"
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "0.4207"
$rng.Style = "Normal"
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "991.33"
$rng.Style = "Normal"
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "7.265"
$rng.Style = "Normal"
$rng = $ws.Range("D49")
$rng.NumberFormat = "@"
$rng.Value = "9.210"
$rng.Style = "Normal"
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "34.10"
$rng.Style = "Normal"

# Column B, C, E cells: plain text assignment (values are not ambiguous numbers).
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -1.98%  "
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("E10").Value = "  +7.19%  "
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("E16").Value = "  +12.66%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("E27").Value = "  +1.48%  "
$ws.Range("E28").Value = "  -2.70%  "
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  -3.41%  "
$ws.Range("E31").Value = "  -6.48%  "
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  -0.52%  "
$ws.Range("E35").Value = "  -2.14%  "
$ws.Range("E36").Value = "  -0.75%  "
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("E39").Value = "  +8.15%  "
$ws.Range("E40").Value = "  -2.88%  "
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("E43").Value = "  +0.90%  "
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("E47").Value = "  +5.23%  "
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("E49").Value = "  +4.88%  "
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").Value = "  -3.19%  "
